# Add a new Q&A row (row 9) to the Sheet1 table for the Postal code analysis.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A9").Value = "Postal code"
$ws.Range("B9").Value = "Get the postal code of Sorgenfrivegen 12B."
$ws.Range("C9").Value = 7031
$ws.Range("D9").Value = "Steps: `n1) Gather and load address point data `n2) Find the ""adresseTekst"" attribute `n3) Find the closest matching address, if any`n4) Get the entire row, and look up the ""postnummer"" attribute"

# Match the author's new cursor/selection position after adding the row.
[void]$ws.Range("C11").Select()
